$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FS-IF-IA")

# Row 24 - Task 1 (B24 = class label "1")
$ws.Range("C24").Value = 0.0
$ws.Range("D24").Value = 0.4
$ws.Range("E24").Value = 0.6
$ws.Range("G24").Value = 57.89473684210527

$ws.Range("J24").Value = 0.9555555555555556
$ws.Range("K24").Value = 0.0
$ws.Range("L24").Value = 0.044444444444444446
$ws.Range("N24").Value = 5.263157894736842

# Row 25 - Task 1 (B25 = class label "2")
$ws.Range("C25").Value = 0.125
$ws.Range("D25").Value = 0.5
$ws.Range("E25").Value = 0.375

$ws.Range("J25").Value = 0.0
$ws.Range("K25").Value = 0.9583333333333334
$ws.Range("L25").Value = 0.041666666666666664

# Row 26 - Task 1 (B26 = class label "3")
$ws.Range("C26").Value = 0.0
$ws.Range("D26").Value = 0.3333333333333333
$ws.Range("E26").Value = 0.6666666666666666

$ws.Range("J26").Value = 0.018518518518518517
$ws.Range("K26").Value = 0.05555555555555555
$ws.Range("L26").Value = 0.9259259259259259
